$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $find"
    }
}

# 1) "a) Przy pożyczkach do kwoty 500,00 zł - sprzedaży bezpośredniej " ->
#    "[#[przedmiot-check-do-500]]Przy pożyczkach do kwoty 500,00 zł - sprzedaży bezpośredniej "
ReplaceText "a) Przy pożyczkach do kwoty 500,00 zł - sprzedaży bezpośredniej " "[#[przedmiot-check-do-500]]Przy pożyczkach do kwoty 500,00 zł - sprzedaży bezpośredniej "

# 2) "b) Przy pożyczkach od kwoty 500,01 zł - aukcji elektronicznej lub sprzedaży bezpośredniej po" ->
#    "[#[przedmiot-check-od-500]]Przy pożyczkach od kwoty 500,01 zł - aukcji elektronicznej lub sprzedaży bezpośredniej po"
ReplaceText "b) Przy pożyczkach od kwoty 500,01 zł - aukcji elektronicznej lub sprzedaży bezpośredniej po" "[#[przedmiot-check-od-500]]Przy pożyczkach od kwoty 500,01 zł - aukcji elektronicznej lub sprzedaży bezpośredniej po"

# 3) "#[przedmiot-procent]%, zostanie zwrócona Pożyczkobiorcy gotówką w terminie 7 dni od dnia otrzymania środków" ->
#    "#[przedmiot-procent]% nadwyżki, zostanie zwrócona Pożyczkobiorcy w terminie 7 dni od dnia otrzymania środków"
ReplaceText "#[przedmiot-procent]%, zostanie zwrócona Pożyczkobiorcy gotówką w terminie 7 dni od dnia otrzymania środków" "#[przedmiot-procent]% nadwyżki, zostanie zwrócona Pożyczkobiorcy w terminie 7 dni od dnia otrzymania środków"

# 4) " przez Pożyczkodawcę w lokalu znajdującym się w miejscowości #[firma-miasto] przy ul. #[firma-adres]." ->
#    " przez Pożyczkodawcę w lokalu znajdującym się w 28-100 Busko-Zdrój ul. Wojska Polskiego 3, bądź przekazem na adres deklarowany w umowie przez pożyczkobiorcę."
ReplaceText " przez Pożyczkodawcę w lokalu znajdującym się w miejscowości #[firma-miasto] przy ul. #[firma-adres]." " przez Pożyczkodawcę w lokalu znajdującym się w 28-100 Busko-Zdrój ul. Wojska Polskiego 3, bądź przekazem na adres deklarowany w umowie przez pożyczkobiorcę."

# 5) "lub na adres siedziby spółki." -> "lub na adres siedziby spółki(28-100 Busko-Zdrój ul. Wojska Polskiego 3)"
ReplaceText "lub na adres siedziby spółki." "lub na adres siedziby spółki(28-100 Busko-Zdrój ul. Wojska Polskiego 3)"

# 6) Remove the "Data przekazania przedmiotu do sprzedaży lub na licytację: #[przedmiot-data-odbioru+30]"
#    paragraph, merging it into the preceding "Maksymalna..." paragraph (only the trailing line break remains).
ReplaceText "naliczonych opłat: #[przedmiot-oplata-max] zł^pData przekazania przedmiotu do sprzedaży lub na licytację: #[przedmiot-data-odbioru+30]" "naliczonych opłat: #[przedmiot-oplata-max] zł"
